# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.382.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.53%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.796.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.59%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.65%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.09"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.28%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.22%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3809"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3472"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.96%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.31"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.12%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.206"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07524"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.72%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.58%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.07"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.21%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.501"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.795.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.065"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.35%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.98%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06663"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.31"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.27%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.527"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.49%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.380.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.74%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.55"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.438"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.35%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.581"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.14%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.34%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.43"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +9.57%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.60"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.96%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.998.16"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.47"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.063"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.96%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.149"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.69%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08711"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.27"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.59%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.701"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.19%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6911"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +10.65%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.450"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.58%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.962"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.17%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06386"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2214"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.80%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02339"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.51%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.275"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.24%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.51"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.84%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6474"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.57%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.869"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.99%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.128"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.47"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.08%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07203"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.53"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.77%  "
